# Apply the "carls newest scraper" update:
#  - drop the columns that are no longer part of the scraped export
#    (Norm/Typ, Ritningsnummer, Position, Beteckning, Kompletterande
#    Information ovrigt, Ref annan, Historiskt Varumarke, Historiskt
#    inkopsreferens, Forpackning, and the always-empty column Q)
#  - keep Varumarke, Artikelbenamning, GVM, Artikelnummer, Typbeteckning,
#    Enhet, SSG-notering, E-nummer, RSK-nummer (now columns A-I)
#  - widen the new SSG-notering column (G) a bit
#  - drop the autofilter / sort state that was scoped to the old layout
#  - refresh the selection / scroll position saved in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete columns, right-to-left so earlier deletes don't
# shift the column letters we still need to remove.
$columnsToRemove = @("Q", "O", "M", "L", "K", "J", "I", "H", "G", "A")
foreach ($col in $columnsToRemove) {
    $ws.Columns($col).Delete() | Out-Null
}

# Give the (new) SSG-notering column a bit more breathing room.
$ws.Columns("G").ColumnWidth = 28.64

# The autofilter / sort state belonged to the old 19-column range; drop it.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Shrink the hidden _FilterDatabase name down to the live data range.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$I`$21217"

# Reset the saved scroll position and refresh the remembered selection to
# match the new (narrower) sheet.
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("H1:H1048576").Select() | Out-Null
